$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data update: insert a new record as the new first data row (row 277),
# shifting all existing records (old rows 277-295) down by one row (to 278-296).
$ws.Rows(277).Insert()

# Populate the newly inserted row 277 with this week's record.
$ws.Range("A277").Value = 4
$ws.Range("B277").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C277").Value = "Los Lagos"
$ws.Range("D277").Value = 44714
$ws.Range("E277").Value = 10
$ws.Range("F277").Value = 100112040
$ws.Range("G277").Value = "Cilantro"
$ws.Range("H277").Value = "Sin especificar"
$ws.Range("I277").Value = "Primera"
$ws.Range("J277").Value = 40
$ws.Range("K277").Value = 6000
$ws.Range("L277").Value = 6000
$ws.Range("M277").Value = 6000
$ws.Range("N277").Value = '$/docena de atados (2 kilos)'
$ws.Range("O277").Value = "Región de La Araucanía"
$ws.Range("P277").Value = 3000
$ws.Range("Q277").Value = 2
$ws.Range("R277").Value = "Hortaliza"
